$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 65: date, hours, note
# Copy the date cell style from the row above (A64) so the new date cell
# reuses the existing date-formatted style instead of creating a new one.
$ws.Range("A64").Copy()
$ws.Range("A65").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A65").Value = 45259
$ws.Range("B65").Value = 3
$ws.Range("C65").Value = "the visa requirements frontend implementation is giving some errors but im trying to fix it"

# Update selection to match new last cell
$ws.Range("C65").Select()
